$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper scratch cell (far away, cleared afterwards) used to coerce a
# numeric-looking value into a genuine text cell (shared string) without
# leaving a new number-format style behind on the target cell.
$scratch = $ws.Cells.Item(1000, 1)

function Set-TextValue($cell, $text) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

function Set-DateValue($cell, $value) {
    $ws.Cells.Item(127, 1).Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats
    $cell.Value = $value
}

# ---- Row 128 ----
Set-DateValue $ws.Cells.Item(128, 1) 45505.2916666667
$ws.Cells.Item(128, 2).Value = 0
$ws.Cells.Item(128, 3).Value = 0.654999971389771
$ws.Cells.Item(128, 4).Value = 0.654999971389771
$ws.Cells.Item(128, 5).Value = 0.654999971389771
$ws.Cells.Item(128, 6).Value = 0.654999971389771
Set-TextValue $ws.Cells.Item(128, 7) "0.654999971389771"
$ws.Cells.Item(128, 8).Value = "BWZ.MI"

# ---- Row 129 ----
Set-DateValue $ws.Cells.Item(129, 1) 45506.6081944444
$ws.Cells.Item(129, 2).Value = 23543
$ws.Cells.Item(129, 3).Value = 0.689999997615814
$ws.Cells.Item(129, 4).Value = 0.639999985694885
$ws.Cells.Item(129, 5).Value = 0.644999980926514
$ws.Cells.Item(129, 6).Value = 0.680000007152557
Set-TextValue $ws.Cells.Item(129, 7) "0.680000007152557"
$ws.Cells.Item(129, 8).Value = "BWZ.MI"
